# Apply the "Updated cryptos list" refresh (prices + 1h volume % deltas),
# and re-rank NEARProtocol/OKB (rows 34 and 35 swap places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the whole Price column to Text first. Price values are displayed
# as literal strings (e.g. "1.00", "0.150") and must not be re-interpreted
# as IEEE-754 numbers (which would lose trailing zeros / introduce rounding
# noise such as 546.03 -> 546.02999999999997).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '62.064.12'
$ws.Range("E2").Value = '  -3.03%  '
$ws.Range("D3").Value = '3.027.20'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '546.03'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").Value = '135.98'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.021.08'
$ws.Range("E8").Value = '  -2.07%  '
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +0.41%  '
$ws.Range("D10").Value = '0.150'
$ws.Range("E10").Value = '  -4.65%  '
$ws.Range("D11").Value = '6.16'
$ws.Range("E11").Value = '  -4.62%  '
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = '34.59'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '3.512.96'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '62.119.96'
$ws.Range("E16").Value = '  -2.95%  '
$ws.Range("E17").Value = '  -2.10%  '
$ws.Range("D18").Value = '3.023.93'
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("D19").Value = '6.71'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '478.69'
$ws.Range("E20").Value = '  -1.31%  '
$ws.Range("D21").Value = '13.31'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '0.679'
$ws.Range("E22").Value = '  -3.24%  '
$ws.Range("D23").Value = '7.12'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '80.75'
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").Value = '12.23'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '2.73'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '1.92'
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").Value = '25.92'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("D33").Value = '2.33'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D36").Value = '5.97'
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("D37").Value = '462.14'
$ws.Range("E37").Value = '  -8.12%  '
$ws.Range("D38").Value = '3.222.77'
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").Value = '0.0802'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("D40").Value = '0.0388'
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("D41").Value = '0.120'
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").Value = '8.19'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").Value = '2.47'
$ws.Range("E43").Value = '  -7.64%  '
$ws.Range("D45").Value = '25.88'
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("D46").Value = '0.246'
$ws.Range("E46").Value = '  -3.79%  '
$ws.Range("D47").Value = '2.02'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("D48").Value = '0.109'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").Value = '118.67'
$ws.Range("E49").Value = '  -3.37%  '
$ws.Range("D50").Value = '0.0₃0498'
$ws.Range("E50").Value = '  -6.64%  '
$ws.Range("E51").Value = '  +6.92%  '

# Rows 34/35 swap: NEARProtocol overtakes OKB in the ranking.
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '5.48'
$ws.Range("E34").Value = '  +2.08%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '55.45'
$ws.Range("E35").Value = '  -3.74%  '

# Restore the default (un-styled) cell style on the Price column; only the
# number format needed forcing to Text, the visual style is unchanged.
$ws.Range("D2:D51").Style = "Normal"
